$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of activity log data
$ws.Range("A5").Value = "Friday 10:00PM - 12:00 AM"
$ws.Range("B5").Value = "Understanding the program"
$ws.Range("C5").Value = "Fix and polish certain areas of the test cases"
$ws.Range("D5").Value = "Rami"

# Adjust column widths to fit new content (values chosen so the engine's
# pixel-snapped stored width lands on the target: 31, 60, 47.19921875, 34.1328125)
$ws.Columns.Item(1).ColumnWidth = 30.166666666666664
$ws.Columns.Item(2).ColumnWidth = 59.16666666666667
$ws.Columns.Item(3).ColumnWidth = 46.33333333333333
$ws.Columns.Item(4).ColumnWidth = 33.33333333333333

# Update the active selection
$ws.Range("C6").Select() | Out-Null
